# Performance.xlsx update:
#  - completed transition from vector to array (pool array replacing pool vector)
#  - added abstract engine class that includes common features
# This adds a new shared string, updates a couple of existing cells/styles,
# refreshes one measured data point (and its dependent formulas), and appends
# a brand-new "latest results" block (rows 71-74) mirroring the layout that is
# already used elsewhere in the sheet (e.g. rows 26-29) for highlighted runs.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) P67 / P69: restyle (re-apply the "note" style already used nearby)
# ---------------------------------------------------------------------------
$ws1.Range("C60").Copy()
$ws1.Range("P67").PasteSpecial(-4122)
$ws1.Range("P67").Value2 = "replaced pool vector with pool array (no board and opposite caching)"

$ws1.Range("C60").Copy()
$ws1.Range("P69").PasteSpecial(-4122)
$ws1.Range("P69").Value2 = "(5932 with board and opposite caching)"

# ---------------------------------------------------------------------------
# 2) E68: updated measurement (dependent formulas F68/G68/H68 recalc automatically)
# ---------------------------------------------------------------------------
$ws1.Range("E68").Value2 = 7448

# ---------------------------------------------------------------------------
# 3) New block: rows 71-74 ("moves cache array" / pool-array results)
#    Formats copied from the analogous highlighted block at rows 26-29.
# ---------------------------------------------------------------------------
$ws1.Range("A26:N27").Copy()
$ws1.Range("A71").PasteSpecial(-4122)

$ws1.Range("A28:N29").Copy()
$ws1.Range("A73").PasteSpecial(-4122)

$ws1.Range("P26").Copy()
$ws1.Range("P71").PasteSpecial(-4122)

# row 71 (depth 4)
$ws1.Range("A71").Value2 = 45873
$ws1.Range("C71").Value2 = 4
$ws1.Range("D71").Value2 = 206603
$ws1.Range("E71").Value2 = 280
$ws1.Range("F71").Formula = "=D71/E71*1000"
$ws1.Range("G71").Formula = "=(E67-E71)/E67"
$ws1.Range("H71").Formula = "=(F71-80000000)/80000000"
$ws1.Range("I71").Value2 = 4
$ws1.Range("J71").Value2 = 197281
$ws1.Range("K71").Value2 = 13
$ws1.Range("L71").Formula = "=J71/K71*1000"
$ws1.Range("M71").Formula = "=(K67-K71)/K67"
$ws1.Range("N71").Formula = "=(L71-80000000)/80000000"
$ws1.Range("P71").Value2 = "moves cache array"

# row 72 (depth 5)
$ws1.Range("C72").Value2 = 5
$ws1.Range("D72").Value2 = 5072212
$ws1.Range("E72").Value2 = 6070
$ws1.Range("F72").Formula = "=D72/E72*1000"
$ws1.Range("G72").Formula = "=(E68-E72)/E68"
$ws1.Range("H72").Formula = "=(F72-80000000)/80000000"
$ws1.Range("I72").Value2 = 5
$ws1.Range("J72").Value2 = 4880523
$ws1.Range("K72").Value2 = 229
$ws1.Range("L72").Formula = "=J72/K72*1000"
$ws1.Range("M72").Formula = "=(K68-K72)/K68"
$ws1.Range("N72").Formula = "=(L72-80000000)/80000000"

# row 73 (depth 6)
$ws1.Range("I73").Value2 = 6
$ws1.Range("J73").Value2 = 119060324
$ws1.Range("K73").Value2 = 4920
$ws1.Range("L73").Formula = "=J73/K73*1000"
$ws1.Range("M73").Formula = "=(K69-K73)/K69"
$ws1.Range("N73").Formula = "=(L73-80000000)/80000000"

# row 74 (depth 7)
$ws1.Range("I74").Value2 = 7
$ws1.Range("J74").Value2 = 3195901860
$ws1.Range("K74").Value2 = 110737
$ws1.Range("L74").Formula = "=J74/K74*1000"
$ws1.Range("N74").Formula = "=(L74-80000000)/80000000"

# Style corrections that deviate from the straight copy-down (match source workbook)
$ws1.Range("H71").Copy()
$ws1.Range("M71").PasteSpecial(-4122)
$ws1.Range("M71").Formula = "=(K67-K71)/K67"

$ws1.Range("H72").Copy()
$ws1.Range("M72").PasteSpecial(-4122)
$ws1.Range("M72").Formula = "=(K68-K72)/K68"

# ---------------------------------------------------------------------------
# 4) Selection/view housekeeping to mirror where the author ended up editing
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("K75").Select()
